# Append two new rows (dates 45985 / 2025-11-24 and 45986 / 2025-11-25) with
# their remn_amt values to the bottom of every worksheet in the workbook.

$wb = $excel.ActiveWorkbook

$newRows = @(
    @{ Date = 45985; Values = @(488, 3362, 2892, 1119) },
    @{ Date = 45986; Values = @(481, 3404, 2814, 1127) }
)

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $lastRow = $ws.Cells.Item(1, 1).End(4).Row   # xlDown = 4, from A1 to last contiguous row
    foreach ($entry in $newRows) {
        $lastRow = $lastRow + 1
        $dateCell = $ws.Cells.Item($lastRow, 1)
        $dateCell.Value = $entry.Date
        $dateCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
        $ws.Cells.Item($lastRow, 2).Value = $entry.Values[$i - 1]
    }
}
